# ETF NAV RTD from 체크 터미널
# Adds a "체크 터미널" style lookup panel (F4:H13) driven by RTD() quotes for
# KODEX 200 (069500), keyed off the ticker entered in G5, alongside the
# pre-existing NAV / 현재가 / 차이 block in columns A:B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Comma-style number format already used by the A1:A4 RTD cells (xfId=1,
# numFmtId 190 -> "_-* #,##0.00_-;\-* #,##0.00_-;_-* "-"_-;_-@_-").
$navFmt = '_-* #,##0.00_-;\-* #,##0.00_-;_-* "-"_-;_-@_-'

# --- New labelled panel: F4/G4 header -----------------------------------
$ws.Range("F4").Value = "종목명"
$ws.Range("G4").Value = "KODEX 200"

# --- Ticker code cell (G5), stored as text, centered --------------------
$ws.Range("F5").Value = "종목코드"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").HorizontalAlignment = -4108
$ws.Range("G5").Value = "069500"

# --- 현재가 (current price) ----------------------------------------------
$ws.Range("F7").Value = "현재가"
$ws.Range("G7").Formula = '=RTD("checkexpert.rtd",,"15001",G5)'
$ws.Range("G7").NumberFormat = $navFmt

# --- 거래량 (volume) ------------------------------------------------------
$ws.Range("F8").Value = "거래량"
$ws.Range("G8").Formula = '=RTD("checkexpert.rtd",,"30620",G5)'
$ws.Range("G8").NumberFormat = $navFmt

# --- iNAV -----------------------------------------------------------------
$ws.Range("F9").Value = "iNAV"
$ws.Range("G9").Formula = '=RTD("checkexpert.rtd",,"15301",G5)'
$ws.Range("G9").NumberFormat = $navFmt

# --- ETF-iNAV ---------------------------------------------------------------
$ws.Range("F10").Value = "ETF-iNAV"
$ws.Range("G10").Formula = "=G7-G9"
$ws.Range("G10").NumberFormat = $navFmt

# --- 괴리율 (%) / tracking error --------------------------------------------
$ws.Range("F11").Value = "괴리율 (%)"
$ws.Range("G11").Formula = "=(G10)/G9*100"
$ws.Range("G11").NumberFormat = $navFmt
$ws.Range("H11").Formula = '=RTD("checkexpert.rtd",,"15304",G5)'
$ws.Range("H11").NumberFormat = "0.00"

# --- 지수기준가 (index reference price) -------------------------------------
$ws.Range("F12").Value = "지수기준가"
$ws.Range("G12").Formula = '=RTD("checkexpert.rtd",,"33405",G5)'
$ws.Range("G12").NumberFormat = $navFmt
$ws.Range("H12").NumberFormat = "0.00"

# --- NAV-지수기준가 ----------------------------------------------------------
$ws.Range("F13").Value = "NAV-지수기준가"
$ws.Range("G13").Formula = "=G9-G12"
$ws.Range("G13").NumberFormat = $navFmt
$ws.Range("H13").Formula = '=RTD("checkexpert.rtd",,"33406",G5)'
$ws.Range("H13").NumberFormat = "0.00"

# --- Column widths to fit the new labels / values ------------------------
$ws.Columns("F").ColumnWidth = 16.125
$ws.Columns("G").ColumnWidth = 13.5

# --- Selection moves to C8, matching the saved view -----------------------
$ws.Range("C8").Select()
